$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D9","D10","D20","D21","D22","D23","D24","D25","D27","D30","D34","D37","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.577.15'
$ws.Range("E2").Value = '  +3.54%  '
$ws.Range("D3").Value = '3.502.64'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '590.82'
$ws.Range("E5").Value = '  +2.93%  '
$ws.Range("D6").Value = '169.56'
$ws.Range("E6").Value = '  +2.80%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '3.499.49'
$ws.Range("E8").Value = '  +1.62%  '
$ws.Range("D9").Value = '0.593'
$ws.Range("E9").Value = '  +6.66%  '
$ws.Range("D10").Value = '7.33'
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("E11").Value = '  +5.60%  '
$ws.Range("E12").Value = '  +3.10%  '
$ws.Range("D13").Value = '4.109.22'
$ws.Range("E13").Value = '  +1.77%  '
$ws.Range("E14").Value = '  -0.84%  '
$ws.Range("E15").Value = '  +3.56%  '
$ws.Range("E16").Value = '  +2.18%  '
$ws.Range("D17").Value = '66.599.39'
$ws.Range("E17").Value = '  +3.46%  '
$ws.Range("D18").Value = '3.517.33'
$ws.Range("E18").Value = '  +2.80%  '
$ws.Range("E19").Value = '  +2.91%  '
$ws.Range("D20").Value = '14.13'
$ws.Range("E20").Value = '  +3.25%  '
$ws.Range("D21").Value = '392.39'
$ws.Range("E21").Value = '  +3.44%  '
$ws.Range("D22").Value = '7.99'
$ws.Range("E22").Value = '  +1.92%  '
$ws.Range("D23").Value = '73.14'
$ws.Range("E23").Value = '  +2.20%  '
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = '0.538'
$ws.Range("E25").Value = '  +3.36%  '
$ws.Range("E26").Value = '  +4.69%  '
$ws.Range("D27").Value = '10.52'
$ws.Range("E27").Value = '  +9.67%  '
$ws.Range("E28").Value = '  +2.13%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = '6.34'
$ws.Range("E30").Value = '  +3.76%  '
$ws.Range("E31").Value = '  +4.23%  '
$ws.Range("E32").Value = '  +2.77%  '
$ws.Range("E33").Value = '  +2.68%  '
$ws.Range("D34").Value = '7.43'
$ws.Range("E34").Value = '  +3.66%  '
$ws.Range("E36").Value = '  +7.32%  '
$ws.Range("D37").Value = '162.55'
$ws.Range("E37").Value = '  +1.56%  '
$ws.Range("E38").Value = '  +2.41%  '
$ws.Range("E39").Value = '  +4.19%  '
$ws.Range("D40").Value = '27.80'
$ws.Range("E40").Value = '  +4.32%  '
$ws.Range("D41").Value = '6.80'
$ws.Range("E41").Value = '  +5.24%  '
$ws.Range("D42").Value = '4.69'
$ws.Range("E42").Value = '  +5.70%  '
$ws.Range("D43").Value = '0.0748'
$ws.Range("E43").Value = '  +2.07%  '
$ws.Range("D44").Value = '26.55'
$ws.Range("E44").Value = '  +1.44%  '
$ws.Range("D45").Value = '2.793.95'
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("D46").Value = '43.25'
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("D47").Value = '0.0311'
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("D48").Value = '2.52'
$ws.Range("E48").Value = '  +0.52%  '
$ws.Range("D49").Value = '350.28'
$ws.Range("E49").Value = '  +4.03%  '
$ws.Range("D50").Value = '1.10'
$ws.Range("E50").Value = '  +3.42%  '
$ws.Range("D51").Value = '33.66'
$ws.Range("E51").Value = '  +11.09%  '
